$wb = $excel.ActiveWorkbook

# --- AntennaMetadata (sheet1) updates ---
$ws = $wb.Worksheets.Item("AntennaMetadata")

# Header row: add SiteCode/River header cells (C1 already = "SiteCode"; add F1 = "River")
$ws.Range("F1").Value = "River"

# Column F (River) - write first occurrences in order River, Colorado River, Fraser River, then
# Connectivity Channel, so shared-string table order matches target.
$ws.Range("F2").Value = "Colorado River"
$ws.Range("F3").Value = "Fraser River"
$ws.Range("F4").Value = "Colorado River"
$ws.Range("F5").Value = "Fraser River"
$ws.Range("F6").Value = "Colorado River"
$ws.Range("F7").Value = "Colorado River"
$ws.Range("F8").Value = "Colorado River"
$ws.Range("F9").Value = "Colorado River"
$ws.Range("F10").Value = "Colorado River"
$ws.Range("F11").Value = "Colorado River"
$ws.Range("F12").Value = "Connectivity Channel"
$ws.Range("F13").Value = "Connectivity Channel"
$ws.Range("F14").Value = "Connectivity Channel"
$ws.Range("F15").Value = "Connectivity Channel"
$ws.Range("F16").Value = "Connectivity Channel"
$ws.Range("F17").Value = "Connectivity Channel"

# --- Notes sheet: insert right after AntennaMetadata ---
$notes = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$notes.Name = "Notes"
$notes.Range("A1").Value = "Assumes that for AntennaMetadata, SiteName and River will not be changing. If they do change, you'll have to go into the runscript and ac couple functions including ""PrepareforMovementsStatesand Summaries"" to change how those variabels are located"
$notes.Range("A2").Select()

# --- Back to AntennaMetadata: update column B (site/antenna names) ---
$ws.Range("B2").Value = "Windy Gap Dam Biomark Antenna"
$ws.Range("B3").Value = "Kaibab Park Biomark Antenna"
$ws.Range("B4").Value = "River Run Biomark Antenna"
$ws.Range("B5").Value = "Fraser River Canyon Biomark Antenna"
$ws.Range("B6").Value = "Red Barn Stationary Antenna"
$ws.Range("B7").Value = "Red Barn Stationary Antenna"
$ws.Range("B8").Value = "Hitching Post Stationary Antenna"
$ws.Range("B9").Value = "Hitching Post Stationary Antenna"
$ws.Range("B10").Value = "Confluence Stationary Antenna"
$ws.Range("B11").Value = "Confluence Stationary Antenna"
$ws.Range("B12").Value = "Connectivity Channel Downstream Stationary Antenna"
$ws.Range("B13").Value = "Connectivity Channel Downstream Stationary Antenna"
$ws.Range("B14").Value = "Connectivity Channel Side Channel Stationary Antenna"
$ws.Range("B15").Value = "Connectivity Channel Side Channel Stationary Antenna"
$ws.Range("B16").Value = "Connectivity Channel Upstream Stationary Antenna"
$ws.Range("B17").Value = "Connectivity Channel Upstream Stationary Antenna"

# New rows 18-19 for Mobile Run tagging stations.
$ws.Range("C18").Value = "M1"
$ws.Range("C19").Value = "M2"
$ws.Range("B18").Value = "Mobile Run"
$ws.Range("B19").Value = "Mobile Run"

$ws.Range("B18").Select()

# --- ImportantStationingVariables sheet: change selection only ---
$isv = $wb.Worksheets.Item("ImportantStationingVariables")
$isv.Range("A3").Select()
$ws.Activate()
